$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "SAN DIEGO AREA TOTALS" label (together with its formatting)
# from B2 to A2
$ws.Range("B2").Copy($ws.Range("A2"))

# Replace B2 with a new "Totals" label, using the worksheet's default
# (unformatted) style
$ws.Range("B2").Value = "Totals"
$ws.Range("B2").Style = "Normal"

# Update the active selection to match the saved view state
$ws.Range("B3").Select()
